# Add a new stock row beneath the existing header row: QR Code (A2),
# Product Name (B2), Quantity (D2) and Selling Price (F2) are filled in;
# Category/Unit/Supplier Name/Warehouse Location (C2/E2/G2/H2) are left
# blank, same as a partially-filled CSV import.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row.
$ws.Range("A2").Value = 133123123123
$ws.Range("B2").Value = "Chungwa"
$ws.Range("D2").Value = 6
$ws.Range("F2").Value = 400

# Touch formatting on the header row and on the new row's "QR Code"/
# "Product Name" cells so they pick up an explicitly-applied cell style
# (distinct from the sheet's implicit default style that the untouched
# numeric cells keep) - mirrors the author re-saving the sheet after
# keying in this row.
$ws.Range("A1:H1").IndentLevel = 1
$ws.Range("A2:B2").IndentLevel = 1

# Leave the active selection on H2, same relative corner as before.
$ws.Range("H2").Select()
